# Update the MU MVD address block: the old post code + "ГСП-1" suffix
# is replaced with the new post code.
#   "...Москва,119991, ГСП-1"  ->  "...Москва,119435"
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "Москва,119991, ГСП-1",  # FindText
    $true,                   # MatchCase
    $false,                  # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    "Москва,119435",         # ReplaceWith
    2                        # Replace (wdReplaceAll)
)

Write-Output "address-zip-updated:$found"
